$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells E/F/G hold percentage/price figures stored as text (not numbers),
# so force a text number format on each target cell before writing the
# new value, to keep Excel from auto-converting them back into numbers.
$textCells = "E2","F2","E3","F3","G3","E4","F4","E5","F5","G5","E6","F6","G6","F7","E8","F8","G8","E9","F9","G9","E10","F10","E11","F11","E12","F12"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("E2").Value = "2.43"
$ws.Range("F2").Value = "-1.73"
$ws.Range("H2").Value = 16811111921.82245

$ws.Range("E3").Value = "2.30"
$ws.Range("F3").Value = "-4.67"
$ws.Range("G3").Value = "0.047138"
$ws.Range("H3").Value = 5910814763.210913

$ws.Range("E4").Value = "-0.01"
$ws.Range("F4").Value = "-0.34"
$ws.Range("H4").Value = 26301615210.77409

$ws.Range("E5").Value = "0.09"
$ws.Range("F5").Value = "0.01"
$ws.Range("G5").Value = "0.009406"
$ws.Range("H5").Value = 396894481.6348014

$ws.Range("E6").Value = "4.12"
$ws.Range("F6").Value = "-4.84"
$ws.Range("G6").Value = "0.002342"
$ws.Range("H6").Value = 1028603983.755987

$ws.Range("F7").Value = "0.04"
$ws.Range("H7").Value = 1947232075.422707

$ws.Range("E8").Value = "2.30"
$ws.Range("F8").Value = "-4.79"
$ws.Range("G8").Value = "0.047105"
$ws.Range("H8").Value = 31166623.05574772

$ws.Range("E9").Value = "0.75"
$ws.Range("F9").Value = "17.24"
$ws.Range("G9").Value = "0.000111"
$ws.Range("H9").Value = 257667042.5942079

$ws.Range("E10").Value = "2.69"
$ws.Range("F10").Value = "-4.89"
$ws.Range("H10").Value = 476588212.4435817

$ws.Range("E11").Value = "8.34"
$ws.Range("F11").Value = "-5.82"
$ws.Range("H11").Value = 546804803.228668

$ws.Range("E12").Value = "1.38"
$ws.Range("F12").Value = "-4.82"
$ws.Range("H12").Value = 150783742.4091178
